# ---------------------------------------------------------------------------
# Adds a new "2022-Q4" quarterly sheet (positioned right after "总计" and
# before "2022-Q3") to the 300026-红日药业 workbook, and updates the "总计"
# (totals) summary sheet so that its first data row reflects the new
# 2022-Q4 figures while every older quarter shifts down by one row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===========================================================================
# STEP 1: update the "总计" (totals) sheet - it is always worksheet #1.
# Row2 becomes the new 2022-Q4 entry, and the old rows 2-5 (2022-Q3 .. 2021-Q4)
# shift down into rows 3-6, with column A renumbered 0..4.
# ===========================================================================
$totalsWs = $wb.Worksheets.Item(1)

$totalsRows = @(
    @(0, '2022-Q4', 18, 7.6),
    @(1, '2022-Q3', 9, 2.56),
    @(2, '2022-Q2', 22, 7.01),
    @(3, '2022-Q1', 10, 1.61),
    @(4, '2021-Q4', 8, 0.71)
)

for ($i = 0; $i -lt $totalsRows.Length; $i++) {
    $r = $i + 2
    $row = $totalsRows[$i]
    $totalsWs.Cells.Item($r, 1).Value = $row[0]
    $totalsWs.Cells.Item($r, 2).Value = $row[1]
    $totalsWs.Cells.Item($r, 3).Value = $row[2]
    $totalsWs.Cells.Item($r, 4).Value = $row[3]
}

# ===========================================================================
# STEP 2: create the new quarterly sheet.
# Worksheets.Add() always inserts the new (blank) sheet at position 1, so
# right after adding, the layout is:
#   1:new  2:总计  3:2022-Q3  4:2022-Q2  5:2022-Q1  6:2021-Q4
# Moving worksheet #1 to sit right before worksheet #3 ("2022-Q3") yields
# the desired final order:
#   1:总计  2:2022-Q4(new)  3:2022-Q3  4:2022-Q2  5:2022-Q1  6:2021-Q4
# NOTE: worksheet object references become stale (they rebind by collection
# index) once the sheets collection is mutated by Add()/Move(), so every
# worksheet must be re-fetched by its *current* index right before use.
# ===========================================================================
$null = $wb.Worksheets.Add()
$wb.Worksheets.Item(1).Move($wb.Worksheets.Item(3))

# From this point on the collection is stable again, so it is safe to keep a
# reference to the new sheet (now at index 2).
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Match the outline / sheetPr settings used by the workbook's other sheets.
$q4.Outline.SummaryRow = 1
$q4.Outline.SummaryColumn = 1

# Match the page margins used by the workbook's other sheets (0.75/0.75/1/1/0.5/0.5 in).
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# Copy the header-row cell format (bold font + border + centered) and the
# column-A cell format from an already-existing quarterly sheet, so the new
# sheet reuses the exact same style that every other quarterly sheet uses.
# ---------------------------------------------------------------------------
$styleSourceWs = $wb.Worksheets.Item(3)   # "2022-Q3" - an existing quarterly sheet

$styleSourceWs.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$styleSourceWs.Range("A2").Copy()
$q4.Range("A2:A19").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Header row.
# ---------------------------------------------------------------------------
$headers = @('基金代码', '基金名称', '基金规模', '股票总仓位', '仓位占比', '持有市值(亿元)', '仓位排名')
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q4.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# ---------------------------------------------------------------------------
# Data rows (18 funds). Columns B, D, E, F and G hold numeric-looking text
# (fund codes with leading zeros, percentages, etc.) that must stay text, so
# the cells are pre-formatted as text ("@") before the value is written, and
# the style is reset back to the default ("Normal") afterwards so the saved
# file does not carry a stray custom style on those cells (matching the
# target, which leaves them with no explicit style).
# ---------------------------------------------------------------------------
$fundData = @(
    @('0', '000772', '景顺长城中国回报灵活配置混合', '35.27', '93.74', '8.05', '2.8392', '8'),
    @('1', '162607', '景顺长城资源垄断混合（LOF）', '24.91', '93.73', '7.93', '1.9754', '7'),
    @('2', '519002', '华安安信消费混合A', '62.23', '83.24', '2.21', '1.3753', '8'),
    @('3', '005505', '前海开源中药研究精选股票A', '5.93', '89.03', '5.86', '0.3475', '8'),
    @('4', '005506', '前海开源中药研究精选股票C', '5.12', '89.03', '5.86', '0.3000', '8'),
    @('5', '013686', '华安安信消费混合C', '10.44', '83.24', '2.21', '0.2307', '8'),
    @('6', '008405', '华泰紫金泰盈混合C', '3.21', '91.44', '6.79', '0.2180', '2'),
    @('7', '011694', '华泰紫金信息科技主题6个月定期开放混合A', '1.42', '92.35', '6.15', '0.0873', '3'),
    @('8', '519673', '银河康乐股票A', '2.15', '93.79', '3.49', '0.0750', '10'),
    @('9', '011695', '华泰紫金信息科技主题6个月定期开放混合C', '0.64', '92.35', '6.15', '0.0394', '3'),
    @('10', '008404', '华泰紫金泰盈混合A', '0.48', '91.44', '6.79', '0.0326', '2'),
    @('11', '011288', '上银医疗健康混合A', '1.17', '88.10', '2.76', '0.0323', '7'),
    @('12', '002681', '金鹰元和灵活配置混合A', '0.30', '81.19', '3.78', '0.0113', '10'),
    @('13', '013920', '兴华创新医疗6个月持有混合A', '0.22', '92.79', '4.12', '0.0091', '7'),
    @('14', '002682', '金鹰元和灵活配置混合C', '0.23', '81.19', '3.78', '0.0087', '10'),
    @('15', '011289', '上银医疗健康混合C', '0.26', '88.10', '2.76', '0.0072', '7'),
    @('16', '016018', '银河康乐股票C', '0.10', '93.79', '3.49', '0.0035', '10'),
    @('17', '013921', '兴华创新医疗6个月持有混合C', '0.06', '92.79', '4.12', '0.0025', '7'),
)

# Pre-format the numeric-looking text columns as text across the whole block.
$q4.Range("B2:B19").NumberFormat = "@"
$q4.Range("D2:G19").NumberFormat = "@"

for ($i = 0; $i -lt $fundData.Length; $i++) {
    $r = $i + 2
    $row = $fundData[$i]
    $q4.Cells.Item($r, 1).Value = [int]$row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 8).Value = [int]$row[7]
}

# Drop the temporary text-number-format back to the default style so the
# cells end up with no explicit style (matching the target workbook).
$q4.Range("B2:B19").Style = "Normal"
$q4.Range("D2:G19").Style = "Normal"

Write-Output "done"
